# Collapse the long tail of account rows (Excel rows 258-386) down to a
# single surviving row (originally row 325: account 004216657 / JOAO /
# -2841.24). Delete from the bottom up so earlier row numbers stay stable.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 326-386 (everything after the row we keep).
$ws.Range("A326:A386").EntireRow.Delete()

# Remove rows 258-324 (everything before the row we keep); the kept row
# (originally 325, JOAO / 004216657 / -2841.24) shifts up to become row 258.
$ws.Range("A258:A324").EntireRow.Delete()
